# Add newer weight-tracker readings (rows 23-37) to the raw_data sheet,
# and refresh the scatter chart that plots A (datetime) vs C (weight) so it
# picks up the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")
$ws.Activate()

# New rows of (datetime, weight) data appended at the bottom of the table.
$newRows = @(
    @(44061.324305555558, 73.599999999999994),
    @(44061.29583333333,  73.599999999999994),
    @(44060.893055555556, 74.2),
    @(44060.322222222225, 74.5),
    @(44060.279861111114, 74.2),
    @(44059.924305555556, 74.8),
    @(44059.241666666669, 74.5),
    @(44059.238888888889, 74.7),
    @(44058.284722222219, 74.7),
    @(44057.328472222223, 74.7),
    @(44057.327777777777, 74.7),
    @(44057.296527777777, 74.7),
    @(44056.936805555553, 75.599999999999994),
    @(44056.352083333331, 75.099999999999994),
    @(44055.93472222222,  76.900000000000006)
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $dt = $newRows[$i][0]
    $w  = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $dt
    $ws.Cells.Item($r, 2).Value = ($dt - [Math]::Floor($dt))
    $ws.Cells.Item($r, 3).Value = $w
    $ws.Cells.Item($r, 4).Formula = '=IF(B' + $r + '<TIME(12,0,0), "AM", "PM")'
}

# Match formatting of the most recently added previous block (rows 19-22).
$fmtSrc = $ws.Range("A22:D22")
$fmtDst = $ws.Range("A23:D37")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)   # xlPasteFormats

$lastRow = $startRow + $newRows.Count - 1

# Update the chart that plots raw_data columns A and C.
$chart = $ws.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = $ws.Range("A2:A$lastRow")
$series.Values = $ws.Range("C2:C$lastRow")

# Move selection the way the author left it.
$ws.Range("C38").Select()
